$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete old data rows (rows 8-13), keep header + rows 2-7 structure
$ws.Range("A8:J13").EntireRow.Delete()

# N/A placeholder value used for columns B-I on each data row
$na = "N/A"

# Row 2
$ws.Range("A2").Value2 = '2025-11-26T20:06:07.721Z'
$ws.Range("B2:I2").Value2 = $na
$ws.Range("J2").Value2 = '{"Sequence":1,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"e96257f1-3a05-4a2f-9c4e-9dd7f68212d4","EventDtm":"2025-11-26T19:58:16Z","AppDtm":"2025-11-26T19:58:35Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-26T19:58:16Z","DeviceID":"JSGA622180044","GPSLockState":"UNLOCKED","GPSSatelliteCount":10,"GPSLastLock":147,"GPSLatitude":26.764466,"GPSLongitude":80.882968,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"405","MNC":"54","LAC":"30917","CID":"234230120","RSSI":"-65","ExtPower":true,"ExtPowerVoltage":32.9,"BatteryVoltage":8.1,"DeviceTemp":34,"RTDLOn":false,"VerFW":"W0206 1.91","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"TRIU6618581","LastAssetRunState":"Running"},"ReeferData":{"ReeferDataDtm":"2025-11-26T19:58:16Z","AssetType":"Reefer","AssetID":"TRIU6618581","OEM":"CARRIER","TAmb":16.86,"TAmbQ":null,"TUSDA4":-50,"TUSDA4Q":"OOR","PctCO2":0,"PctCO2Q":null,"PctCO2Set":5,"PctCO2SetQ":null,"PSuc":-14.73,"PSucQ":"OOR","TDis":-196.88,"TDisQ":"OOR","FreqComp":null,"TSuc":0,"TSucQ":null,"MCond":1,"PCond":104.95,"PCondQ":null,"TCond":null,"TCondQ":null,"MCtrl":"Cool","SnCtrl":"04900648","AmpPhA":11.75,"AmpPhB":12.34,"AmpPhC":11.75,"PDis":-14.6,"PDisQ":"OOR","PctEconVlv":0,"PctEconVlvQ":null,"PctExpVlv":0,"PctExpVlvQ":null,"TEvap":-1.53,"TEvapQ":null,"MCtrl3":null,"PctHtr":null,"PctHtrQ":null,"MEvapFanHS":1,"PctGasVlv":null,"PctGasVlvQ":null,"PctHum":96.34,"PctHumQ":null,"PctHumSet":0,"PctHumSetQ":"OOR","FreqLine":50,"FreqLineQ":null,"VLine1":426.19,"VLine2":null,"VLine3":null,"MEvapFanLS":0,"PctO2":0,"PctO2Q":"asProvided","PctO2Set":10,"PctO2SetQ":null,"MCtrl2":null,"TRtn1":-2.35,"TRtn1Q":null,"TRtn2":-2.38,"TRtn2Q":null,"TSet":-3,"TSetQ":null,"VerSWMajor":"5167","VerSWMinor":null,"PctSucVlv":10.43,"PctSucVlvQ":null,"TSup1":-3.02,"TSup1Q":null,"TSup2":-2.91,"TSup2Q":null,"AmpTotalDraw":null,"AmpTotalDrawQ":null,"TUSDA1":-50,"TUSDA1Q":"OOR","TUSDA2":-50,"TUSDA2Q":"OOR","TUSDA3":-50,"TUSDA3Q":"OOR","CmhVent":0,"CmhVentQ":null,"BkNum":null,"TBk":null,"TBkQ":null,"PTIDtm":"2025-06-02T11:27:53Z","PTIResult":"Skipped","TWResult":null,"TWExpiration":null,"TWExpirationQ":null,"AtmosMode":null,"ReeferAlarms":[{"OemAlarm":52,"RCAlias":"AL52","RCSeverity":"Minor","Active":true},{"OemAlarm":60,"RCAlias":"AL60","RCSeverity":"Informational","Active":true}]}}}'

# Row 3
$ws.Range("A3").Value2 = '2025-11-26T20:06:07.796Z'
$ws.Range("B3:I3").Value2 = $na
$ws.Range("J3").Value2 = '{"Sequence":2,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"e7c568d3-aaad-49cf-ba88-7205b0f100b1","EventDtm":"2025-11-26T19:59:44Z","AppDtm":"2025-11-26T20:00:03Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-26T19:59:44Z","DeviceID":"JSGA623040321","GPSLockState":"LOCKED","GPSSatelliteCount":10,"GPSLastLock":0,"GPSLatitude":12.816277,"GPSLongitude":77.68725,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"45","LAC":"9057","CID":"180498700","RSSI":"-55","ExtPower":true,"ExtPowerVoltage":31.8,"BatteryVoltage":8,"DeviceTemp":26,"RTDLOn":false,"VerFW":"W0206 1.41","DeviceType":"CT3500","DoorState":"Open","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"TRIU6617250","LastAssetRunState":"Running"},"ReeferData":null}}'

# Row 4
$ws.Range("A4").Value2 = '2025-11-26T20:06:07.831Z'
$ws.Range("B4:I4").Value2 = $na
$ws.Range("J4").Value2 = '{"Sequence":3,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"83568cd2-b74c-42b6-bc65-4b36a010fcce","EventDtm":"2025-11-26T20:03:56Z","AppDtm":"2025-11-26T20:04:10Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-26T20:03:56Z","DeviceID":"JSGA621320217","GPSLockState":"LOCKED","GPSSatelliteCount":14,"GPSLastLock":0,"GPSLatitude":26.644384,"GPSLongitude":88.468068,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"405","MNC":"51","LAC":"41512","CID":"16323842","RSSI":null,"ExtPower":true,"ExtPowerVoltage":33.7,"BatteryVoltage":8.1,"DeviceTemp":24,"RTDLOn":false,"VerFW":"W0206 1.75","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"CXRU1036301","LastAssetRunState":"Running"},"ReeferData":{"ReeferDataDtm":"2025-11-26T20:03:56Z","AssetType":"Reefer","AssetID":"CXRU1036301","OEM":"DAIKIN","TAmb":18.81,"TAmbQ":null,"TUSDA4":-53.4,"TUSDA4Q":"OOR","PctCO2":25.5,"PctCO2Q":"OOR","PctCO2Set":25.5,"PctCO2SetQ":"OOR","PSuc":50,"PSucQ":"asProvided","TDis":77.2,"TDisQ":null,"FreqComp":null,"TSuc":-14.9,"TSucQ":null,"MCond":"On","PCond":null,"PCondQ":null,"TCond":null,"TCondQ":null,"MCtrl":null,"SnCtrl":null,"AmpPhA":null,"AmpPhB":null,"AmpPhC":null,"PDis":810,"PDisQ":"asProvided","PctEconVlv":0,"PctEconVlvQ":null,"PctExpVlv":14,"PctExpVlvQ":null,"TEvap":-0.62,"TEvapQ":null,"MCtrl3":"Modulation","PctHtr":null,"PctHtrQ":null,"MEvapFanHS":"On","PctGasVlv":0,"PctGasVlvQ":null,"PctHum":100.39,"PctHumQ":"OOR","PctHumSet":65,"PctHumSetQ":"off","FreqLine":null,"FreqLineQ":null,"VLine1":423.5,"VLine2":null,"VLine3":null,"MEvapFanLS":"Off","PctO2":25.5,"PctO2Q":"OOR","PctO2Set":25.5,"PctO2SetQ":"OOR","MCtrl2":"Modulation","TRtn1":0.81,"TRtn1Q":null,"TRtn2":1.6,"TRtn2Q":null,"TSet":0,"TSetQ":null,"VerSWMajor":"265A","VerSWMinor":null,"PctSucVlv":17.99,"PctSucVlvQ":null,"TSup1":0.12,"TSup1Q":null,"TSup2":1.8,"TSup2Q":null,"AmpTotalDraw":12,"AmpTotalDrawQ":"asProvided","TUSDA1":-53.4,"TUSDA1Q":"OOR","TUSDA2":-53.4,"TUSDA2Q":"OOR","TUSDA3":-53.4,"TUSDA3Q":"OOR","CmhVent":1020,"CmhVentQ":"error","BkNum":null,"TBk":null,"TBkQ":null,"PTIDtm":"2025-07-16T10:29:07Z","PTIResult":"Failed","TWResult":null,"TWExpiration":null,"TWExpirationQ":null,"AtmosMode":"UnitOff","ReeferAlarms":null}}}'

# Row 5
$ws.Range("A5").Value2 = '2025-11-26T20:06:08.393Z'
$ws.Range("B5:I5").Value2 = $na
$ws.Range("J5").Value2 = '{"Sequence":4,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"925e9264-9281-4fb6-8fdb-2d25a0bbc0ef","EventDtm":"2025-11-26T20:04:16Z","AppDtm":"2025-11-26T20:04:33Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-26T20:04:16Z","DeviceID":"JSGA623040311","GPSLockState":"LOCKED","GPSSatelliteCount":10,"GPSLastLock":0,"GPSLatitude":26.606132,"GPSLongitude":80.723723,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"15","LAC":"55363","CID":"210193710","RSSI":"-79","ExtPower":true,"ExtPowerVoltage":34,"BatteryVoltage":8,"DeviceTemp":14,"RTDLOn":false,"VerFW":"W0206 1.75","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":null,"LastAssetRunState":null},"ReeferData":null}}'

# Row 6
$ws.Range("A6").Value2 = '2025-11-26T20:06:08.953Z'
$ws.Range("B6:I6").Value2 = $na
$ws.Range("J6").Value2 = '{"Sequence":5,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"5805fde2-4ad1-4c70-ae25-620c91b315cf","EventDtm":"2025-11-26T20:04:33Z","AppDtm":"2025-11-26T20:04:49Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-26T20:04:33Z","DeviceID":"JSGA622340343","GPSLockState":"LOCKED","GPSSatelliteCount":13,"GPSLastLock":0,"GPSLatitude":12.805029,"GPSLongitude":77.662178,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"45","LAC":"9003","CID":"46904076","RSSI":null,"ExtPower":true,"ExtPowerVoltage":29.7,"BatteryVoltage":8,"DeviceTemp":24,"RTDLOn":false,"VerFW":"W0206 1.75","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"CCLU1035976","LastAssetRunState":"Running"},"ReeferData":{"ReeferDataDtm":"2025-11-26T20:04:33Z","AssetType":"Reefer","AssetID":"CCLU1035976","OEM":"CARRIER","TAmb":28.31,"TAmbQ":null,"TUSDA4":-50,"TUSDA4Q":"OOR","PctCO2":5,"PctCO2Q":null,"PctCO2Set":1,"PctCO2SetQ":null,"PSuc":-14.73,"PSucQ":"OOR","TDis":-196.88,"TDisQ":"OOR","FreqComp":null,"TSuc":0,"TSucQ":null,"MCond":1,"PCond":149.42,"PCondQ":null,"TCond":null,"TCondQ":null,"MCtrl":"Cool","SnCtrl":"04904656","AmpPhA":10.18,"AmpPhB":11.1,"AmpPhC":10.18,"PDis":-14.6,"PDisQ":"OOR","PctEconVlv":0,"PctEconVlvQ":null,"PctExpVlv":0,"PctExpVlvQ":null,"TEvap":6.06,"TEvapQ":null,"MCtrl3":null,"PctHtr":null,"PctHtrQ":null,"MEvapFanHS":1,"PctGasVlv":null,"PctGasVlvQ":null,"PctHum":0,"PctHumQ":null,"PctHumSet":0,"PctHumSetQ":"OOR","FreqLine":50,"FreqLineQ":null,"VLine1":417.63,"VLine2":null,"VLine3":null,"MEvapFanLS":0,"PctO2":0.11,"PctO2Q":"unknown","PctO2Set":3,"PctO2SetQ":null,"MCtrl2":null,"TRtn1":5.61,"TRtn1Q":null,"TRtn2":5.63,"TRtn2Q":null,"TSet":5,"TSetQ":null,"VerSWMajor":"5180","VerSWMinor":null,"PctSucVlv":3.75,"PctSucVlvQ":null,"TSup1":5.09,"TSup1Q":null,"TSup2":5.03,"TSup2Q":null,"AmpTotalDraw":null,"AmpTotalDrawQ":null,"TUSDA1":-50,"TUSDA1Q":"OOR","TUSDA2":-50,"TUSDA2Q":"OOR","TUSDA3":-50,"TUSDA3Q":"OOR","CmhVent":0,"CmhVentQ":null,"BkNum":null,"TBk":null,"TBkQ":null,"PTIDtm":"2025-05-12T10:01:32Z","PTIResult":"Failed","TWResult":null,"TWExpiration":null,"TWExpirationQ":null,"AtmosMode":null,"ReeferAlarms":null}}}'

# Row 7
$ws.Range("A7").Value2 = '2025-11-26T20:07:04.694Z'
$ws.Range("B7:I7").Value2 = $na
$ws.Range("J7").Value2 = '{"Sequence":6,"Event":{"EventClass":"DeviceMessage","MessageData":{"MsgID":"ba8e354d-f391-48d1-b488-e8709d49fc2c","EventDtm":"2025-11-26T20:06:47Z","AppDtm":"2025-11-26T20:07:04Z","Events":["ScheduledReport"]},"DeviceData":{"DeviceDataDtm":"2025-11-26T20:06:47Z","DeviceID":"JSGA621320194","GPSLockState":"LOCKED","GPSSatelliteCount":15,"GPSLastLock":0,"GPSLatitude":16.824892,"GPSLongitude":78.131558,"GeofenceCode":null,"ServerGeofenceCode":null,"MCC":"404","MNC":"49","LAC":"31204","CID":"255654657","RSSI":null,"ExtPower":true,"ExtPowerVoltage":31.5,"BatteryVoltage":8,"DeviceTemp":25,"RTDLOn":false,"VerFW":"W0206 1.91","DeviceType":"CT3500","DoorState":"Disconnected","DoorStateDtm":null,"DoorSensorBatteryVoltage":null,"LastAssetID":"CXRU1186425","LastAssetRunState":"Running"},"ReeferData":{"ReeferDataDtm":"2025-11-26T20:06:47Z","AssetType":"Reefer","AssetID":"CXRU1186425","OEM":"DAIKIN","TAmb":20.62,"TAmbQ":null,"TUSDA4":-46.5,"TUSDA4Q":"OOR","PctCO2":25.5,"PctCO2Q":"OOR","PctCO2Set":25.5,"PctCO2SetQ":"OOR","PSuc":-20,"PSucQ":"asProvided","TDis":57.6,"TDisQ":null,"FreqComp":null,"TSuc":-28.3,"TSucQ":null,"MCond":"On","PCond":null,"PCondQ":null,"TCond":null,"TCondQ":null,"MCtrl":null,"SnCtrl":null,"AmpPhA":null,"AmpPhB":null,"AmpPhC":null,"PDis":800,"PDisQ":"asProvided","PctEconVlv":0,"PctEconVlvQ":null,"PctExpVlv":37,"PctExpVlvQ":null,"TEvap":9.62,"TEvapQ":null,"MCtrl3":"Modulation","PctHtr":null,"PctHtrQ":null,"MEvapFanHS":"On","PctGasVlv":0,"PctGasVlvQ":null,"PctHum":100.39,"PctHumQ":"OOR","PctHumSet":65,"PctHumSetQ":"off","FreqLine":null,"FreqLineQ":null,"VLine1":395.8,"VLine2":null,"VLine3":null,"MEvapFanLS":"Off","PctO2":25.5,"PctO2Q":"OOR","PctO2Set":25.5,"PctO2SetQ":"OOR","MCtrl2":"Modulation","TRtn1":4.62,"TRtn1Q":null,"TRtn2":4.4,"TRtn2Q":null,"TSet":5,"TSetQ":null,"VerSWMajor":"265A","VerSWMinor":null,"PctSucVlv":9.45,"PctSucVlvQ":null,"TSup1":5.12,"TSup1Q":null,"TSup2":5.5,"TSup2Q":null,"AmpTotalDraw":16,"AmpTotalDrawQ":"asProvided","TUSDA1":-46.5,"TUSDA1Q":"OOR","TUSDA2":-46.5,"TUSDA2Q":"OOR","TUSDA3":-46.5,"TUSDA3Q":"OOR","CmhVent":1020,"CmhVentQ":"error","BkNum":null,"TBk":null,"TBkQ":null,"PTIDtm":"2025-06-30T05:28:13Z","PTIResult":"Passed","TWResult":null,"TWExpiration":null,"TWExpirationQ":null,"AtmosMode":"UnitOff","ReeferAlarms":null}}}'
